$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E25's bsecode was stored as text ("544028"); correct it to a real number.
$ws.Range("E25").Value = 544028

# Append the new screener row (row 26) that was scraped for this breakout run.
$ws.Range("A26").Value = "21/06/2024 06:45:24"
$ws.Range("B26").Value = 1
$ws.Range("C26").Value = "TATATECH"
$ws.Range("D26").Value = "Tata Technologies Ltd"

# bsecode keeps its original (text) representation for this row, so force
# text formatting before assigning the numeric-looking string, then drop
# back to the sheet's default style so no extra formatting is left behind.
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "544028"
$ws.Range("E26").Style = "Normal"

$ws.Range("F26").Value = -0.73
$ws.Range("G26").Value = 1002.95
$ws.Range("H26").Value = 1794692
